$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 16: Registration margin left ---
$ws.Range("A16").Value = "Registration margin left"
$ws.Range("B16").Value = 54
$ws.Range("C16").Value = 910
$ws.Range("D16").Formula = "=B16/C16"
$ws.Range("D16").NumberFormat = "0.00%"
$ws.Range("E16").Value = 411
$ws.Range("E16").NumberFormat = "0"
$ws.Range("F16").Formula = "=D16*E16"
$ws.Range("F16").NumberFormat = "0.00"

# --- Row 17: Registration top margin text field ---
$ws.Range("A17").Value = "Registration top margin text field"
$ws.Range("B17").Value = 26
$ws.Range("C17").Value = 944
$ws.Range("D17").Formula = "=B17/C17"
$ws.Range("D17").NumberFormat = "0.00%"
$ws.Range("E17").Value = 891
$ws.Range("E17").NumberFormat = "0.00"
$ws.Range("F17").Formula = "=D17*E17"
$ws.Range("F17").NumberFormat = "0.00"

# --- Row 18: Registration top margin text view ---
$ws.Range("A18").Value = "Registration top margin text view"
$ws.Range("B18").Value = 50
$ws.Range("C18").Value = 944
$ws.Range("D18").Formula = "=B18/C18"
$ws.Range("D18").NumberFormat = "0.00%"
$ws.Range("E18").Value = 891
$ws.Range("E18").NumberFormat = "0.00"
$ws.Range("F18").Formula = "=D18*E18"
$ws.Range("F18").NumberFormat = "0.00"

# --- Row 19: Registration top margin text view (again) ---
$ws.Range("A19").Value = "Registration top margin text view"
$ws.Range("B19").Value = 17
$ws.Range("C19").Value = 944
$ws.Range("D19").Formula = "=B19/C19"
$ws.Range("D19").NumberFormat = "0.00%"
$ws.Range("E19").Value = 891
$ws.Range("E19").NumberFormat = "0.00"
$ws.Range("F19").Formula = "=D19*E19"
$ws.Range("F19").NumberFormat = "0.00"

# --- Update view/selection: scroll back to top, select A6 ---
$ws.Range("A6").Select() | Out-Null
